# Merge the split "<id>...</id>" runs back into a single run, for each
# of the three occurrences in the document. Previously the id value was
# typed as a separate (differently-formatted) run sandwiched between a
# run holding "<id>" and a run holding "</id>"; now all three collapse
# into one run carrying the original "<id>" run's formatting.

$d = $word.ActiveDocument

$ids = @("p098v_5", "p099r_2", "p099r_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2)
}
